# Generate Report for Archive
# Change the "Ready for handoff" status text to "In Translation" across all
# sheets that reference it (Overview, zh-cn, de-de), then let the column
# widths reflect the new (shorter) text via AutoFit, matching the narrower
# column widths recorded in the target workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: keep the string literal on the LEFT of -eq. PowerShell's
        # -eq coerces the right-hand side to the left operand's type, so
        # "text" -eq $cell.Value() compares as strings; the reverse
        # ($cell.Value() -eq "text") would coerce the literal to the left
        # side's type and wrongly match boolean/numeric cells (e.g. a
        # boolean TRUE cell would coerce "Ready for handoff" to $true).
        if ("Ready for handoff" -eq $cell.Value()) {
            $cell.Value = "In Translation"
        }
    }
}

# Autofit the affected columns on each sheet so widths shrink to fit the
# new (shorter) text, same as Excel does automatically when a report is
# regenerated. Use EntireColumn.AutoFit() off a Range (rather than
# Columns.Item("E:F")) for reliable multi-column addressing.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E1:F1").EntireColumn.AutoFit() | Out-Null

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C1").EntireColumn.AutoFit() | Out-Null

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C1").EntireColumn.AutoFit() | Out-Null
